$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $value) {
    # Prefix with an apostrophe so Excel stores the value as text
    # (matching the original inline-string / text cell type) instead of
    # auto-converting numeric-looking strings into numbers.
    $ws.Range($addr).Value = "'" + $value
}

function Set-PlainCell($addr, $value) {
    $ws.Range($addr).Value = $value
}

# Row 2 - Bitcoin
Set-TextCell "D2" "30.103.02"
Set-PlainCell "E2" "  +2.07%  "

# Row 3 - Ethereum
Set-TextCell "D3" "1.949.75"
Set-PlainCell "E3" "  +1.46%  "

# Row 4 - TetherUSD
Set-PlainCell "E4" "  +0.12%  "

# Row 5 - BNB
Set-TextCell "D5" "327.80"
Set-PlainCell "E5" "  +0.73%  "

# Row 6 - USDC
Set-PlainCell "E6" "  +0.29%  "

# Row 7 - XRP
Set-TextCell "D7" "0.4857"
Set-PlainCell "E7" "  +0.37%  "

# Row 8 - Cardano
Set-TextCell "D8" "0.4118"
Set-PlainCell "E8" "  +0.49%  "

# Row 9 - Dogecoin
Set-TextCell "D9" "0.08239"
Set-PlainCell "E9" "  +0.73%  "

# Row 10 - Polygon
Set-TextCell "D10" "1.021"
Set-PlainCell "E10" "  -0.33%  "

# Row 11 - Solana
Set-TextCell "D11" "24.07"
Set-PlainCell "E11" "  +2.24%  "

# Row 12 - WrappedEther
Set-TextCell "D12" "1.941.49"
Set-PlainCell "E12" "  -0.13%  "

# Row 13 - Polkadot
Set-TextCell "D13" "6.124"
Set-PlainCell "E13" "  +1.25%  "

# Row 14 - Chainlink
Set-TextCell "D14" "7.343"
Set-PlainCell "E14" "  +1.30%  "

# Row 15 - Litecoin
Set-TextCell "D15" "92.05"
Set-PlainCell "E15" "  +0.69%  "

# Row 16 - TRON
Set-TextCell "D16" "0.06867"
Set-PlainCell "E16" "  +1.33%  "

# Row 17 - BinanceUSD
Set-PlainCell "E17" "  +0.28%  "

# Row 18 - ShibaInu
Set-PlainCell "E18" "  +0.14%  "

# Row 19 - Avalanche
Set-TextCell "D19" "17.92"
Set-PlainCell "E19" "  +0.83%  "

# Row 20 - Dai
Set-TextCell "D20" "1.009"
Set-PlainCell "E20" "  +0.33%  "

# Row 21 - WrappedBTC
Set-TextCell "D21" "30.099.18"
Set-PlainCell "E21" "  +1.95%  "

# Row 22 - Uniswap
Set-TextCell "D22" "5.699"

# Row 23 - Cosmos
Set-TextCell "D23" "12.01"
Set-PlainCell "E23" "  +2.00%  "

# Row 24 - Toncoin
Set-TextCell "D24" "2.201"
Set-PlainCell "E24" "  +0.69%  "

# Row 25 - WrappedliquidstakedEther2.0
Set-TextCell "D25" "2.181.28"
Set-PlainCell "E25" "  +2.88%  "

# Row 26 - InternetComputer(DFINITY)
Set-TextCell "D26" "6.571"
Set-PlainCell "E26" "  -2.58%  "

# Row 27 - Monero
Set-TextCell "D27" "157.03"
Set-PlainCell "E27" "  +0.22%  "

# Row 28 - EthereumClassic
Set-TextCell "D28" "20.19"
Set-PlainCell "E28" "  +0.69%  "

# Row 29 - LidoDAOToken
Set-PlainCell "E29" "  -0.07%  "

# Row 30 - BitcoinCash
Set-TextCell "D30" "121.51"
Set-PlainCell "E30" "  +0.75%  "

# Row 31 - ImmutableX
Set-PlainCell "E31" "  -0.41%  "

# Row 32 - Stellar
Set-TextCell "D32" "0.09648"
Set-PlainCell "E32" "  +0.66%  "

# Row 33 - Filecoin
Set-TextCell "D33" "5.655"
Set-PlainCell "E33" "  +2.29%  "

# Row 34 - ARBITRUM
Set-TextCell "D34" "1.432"
Set-PlainCell "E34" "  +2.88%  "

# Row 35 - HuobiToken
Set-PlainCell "E35" "  -0.28%  "

# Row 36 - Hedera
Set-TextCell "D36" "0.06533"
Set-PlainCell "E36" "  +6.27%  "

# Row 37 - VeChain
Set-TextCell "D37" "0.02311"
Set-PlainCell "E37" "  +1.04%  "

# Row 38 - TrustWalletToken
Set-TextCell "D38" "1.237"
Set-PlainCell "E38" "  +4.79%  "

# Row 39 - TheSandbox
Set-TextCell "D39" "0.5993"
Set-PlainCell "E39" "  +0.07%  "

# Row 40 - Aptos
Set-TextCell "D40" "10.80"
Set-PlainCell "E40" "  +0.01%  "

# Row 41 - FraxShare
Set-TextCell "D41" "8.006"
Set-PlainCell "E41" "  -0.45%  "

# Row 42 - RenderToken
Set-TextCell "D42" "2.546"
Set-PlainCell "E42" "  +5.66%  "

# Row 43 - Algorand
Set-TextCell "D43" "0.1861"
Set-PlainCell "E43" "  -0.19%  "

# Row 44 & 45 - swap EnergySwap / WEMIXToken (rank order changes)
Set-PlainCell "B44" "WEMIXToken"
Set-PlainCell "C44" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell "D44" "1.284"
Set-PlainCell "E44" "  +0.14%  "

Set-PlainCell "B45" "EnergySwap"
Set-PlainCell "C45" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell "D45" "12.45"
Set-PlainCell "E45" "  +0.06%  "

# Row 46 - Cronos
Set-TextCell "D46" "0.07565"
Set-PlainCell "E46" "  -0.51%  "

# Row 47 - Decentraland
Set-TextCell "D47" "0.5609"
Set-PlainCell "E47" "  +0.32%  "

# Row 48 - NEARProtocol
Set-TextCell "D48" "1.997"
Set-PlainCell "E48" "  +1.78%  "

# Row 49 - Quant
Set-TextCell "D49" "118.01"
Set-PlainCell "E49" "  +0.60%  "

# Row 50 - MXToken
Set-TextCell "D50" "2.444"
Set-PlainCell "E50" "  +0.21%  "

# Row 51 - Aave
Set-TextCell "D51" "72.75"
Set-PlainCell "E51" "  +0.07%  "
